$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at G (shifts the existing FIYAT column from G to H)
# and resize column E to better fit the now-narrower layout.
$ws.Range("G1").EntireColumn.Insert()
$ws.Range("G1").Value = "KDV"
$ws.Columns.Item(5).ColumnWidth = 57.7

# Select the newly inserted column, matching the author's final UI state.
$ws.Range("G:G").Select()
